# The "Recorded By" column (G) lists who touched each attendance session,
# e.g. "dnasr281@gmail.com, System" or "backup@backdoor.com, system, System".
# This normalizes each list so the automated "System" entry (exact case)
# is always reported first, leaving the human/other recorders after it in
# their original relative order. Cells without an exact "System" entry
# (including ones that are only lowercase "system", or have no System at
# all) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dims = $ws.UsedRange
$lastRow = $dims.Rows.Count + $dims.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G
    $val = $cell.Value2
    if ($val -eq $null) { continue }
    if ($val -notlike "*System*") { continue }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -eq 0) { continue }

    $newParts = $systemParts + $otherParts
    $newVal = $newParts -join ", "

    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
